# smooth interpolation of SIMTRA profiles. analytycal integration of 1-st order interpolants
#
# Adds a new "Сглаживание профиля" (smooth) parameter row to the SIMTRA
# parameter sheet, and restores the view/selection state recorded for the
# sheet after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append the new parameter row (row 46) -------------------------------
$ws.Range("A46").Value = 44
$ws.Range("B46").Value = "Сглаживание профиля"
$ws.Range("C46").Value = "smooth"
$ws.Range("D46").Value = 0
$ws.Range("F46").Value = "0+float"
$ws.Range("G46").Value = "model"

# --- Restore the sheet's active selection --------------------------------
$ws.Range("J49").Select()
